$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.810.74'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.463.02'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '574.55'
$ws.Range('D6').Value = '146.89'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = '2.463.06'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('D11').Value = '0.163'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '28.98'
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '2.910.42'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').Value = '62.734.86'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').Value = '2.459.69'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').Value = '7.94'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '326.73'
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('E23').Value = '  +8.86%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '10.06'
$ws.Range('E25').Value = '  +18.31%  '
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').Value = '646.89'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '0.0₃0984'
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -15.29%  '
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  +2.86%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = '2.81'
$ws.Range('E38').Value = '  +3.68%  '
$ws.Range('D39').Value = '152.16'
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -2.06%  '
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').Value = '0.0₆0308'
$ws.Range('E44').Value = '  -39.81%  '
$ws.Range('D46').Value = '152.23'
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('D47').Value = '15.27'
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('D48').Value = '3.59'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('D49').Value = '20.51'
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('E51').Value = '  -1.12%  '
